$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, pushing the existing row 4 (and below) down
$ws.Rows.Item(4).Insert()

# Match the custom row height used by all the other data rows
$ws.Rows.Item(4).RowHeight = 18.75

# Fill the new row 4 with the "html css" action entry
$ws.Cells.Item(4, 1).Value = "html css"
$ws.Cells.Item(4, 2).Value = "_"
$ws.Cells.Item(4, 3).Value = "0h 5m"
$ws.Cells.Item(4, 4).Value = 1

# Re-select cell B5 as shown in the diff (activeCell moved there after the edit)
$ws.Range("B5").Select()
